# Apply updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell while forcing text interpretation
# (prefixing with a literal leading apostrophe), so numeric-looking
# strings such as "4.00" or "0.140" are not auto-converted to numbers
# and keep their original formatting/precision, matching the source data.
function Set-TextValue($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'$text"
}

# Row 2
Set-TextValue 2 4 '62.085.88'
$ws.Cells.Item(2, 5).Value = '  +1.73%  '

# Row 3
Set-TextValue 3 4 '3.435.03'
$ws.Cells.Item(3, 5).Value = '  +1.96%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.08%  '

# Row 5
Set-TextValue 5 4 '409.22'
$ws.Cells.Item(5, 5).Value = '  +0.77%  '

# Row 6
Set-TextValue 6 4 '129.54'
$ws.Cells.Item(6, 5).Value = '  -3.03%  '

# Row 7
Set-TextValue 7 4 '0.631'
$ws.Cells.Item(7, 5).Value = '  +6.68%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.12%  '

# Row 9
Set-TextValue 9 4 '0.759'
$ws.Cells.Item(9, 5).Value = '  +13.04%  '

# Row 10
Set-TextValue 10 4 '0.140'
$ws.Cells.Item(10, 5).Value = '  +15.54%  '

# Row 11
Set-TextValue 11 4 '43.31'
$ws.Cells.Item(11, 5).Value = '  +2.18%  '

# Row 12
$ws.Cells.Item(12, 5).Value = '  -0.45%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +5.87%  '

# Row 14
Set-TextValue 14 4 '20.31'
$ws.Cells.Item(14, 5).Value = '  +3.64%  '

# Row 15
Set-TextValue 15 4 '0.0000193'
$ws.Cells.Item(15, 5).Value = '  +51.37%  '

# Row 16
Set-TextValue 16 4 '3.432.43'
$ws.Cells.Item(16, 5).Value = '  +2.26%  '

# Row 17
$ws.Cells.Item(17, 2).Value = 'WrappedBTC'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 17 4 '62.149.34'
$ws.Cells.Item(17, 5).Value = '  +1.85%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'Polygon'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 18 4 '1.05'
$ws.Cells.Item(18, 5).Value = '  +2.90%  '

# Row 19
Set-TextValue 19 4 '11.39'
$ws.Cells.Item(19, 5).Value = '  +3.03%  '

# Row 20
Set-TextValue 20 4 '376.19'
$ws.Cells.Item(20, 5).Value = '  +22.23%  '

# Row 21
Set-TextValue 21 4 '87.70'
$ws.Cells.Item(21, 5).Value = '  +4.28%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  -1.15%  '

# Row 23
Set-TextValue 23 4 '13.37'
$ws.Cells.Item(23, 5).Value = '  +5.28%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +2.82%  '

# Row 25
Set-TextValue 25 4 '31.61'
$ws.Cells.Item(25, 5).Value = '  +7.54%  '

# Row 26
Set-TextValue 26 4 '4.81'
$ws.Cells.Item(26, 5).Value = '  +0.63%  '

# Row 27
Set-TextValue 27 4 '8.46'
$ws.Cells.Item(27, 5).Value = '  +2.04%  '

# Row 29
Set-TextValue 29 4 '2.74'
$ws.Cells.Item(29, 5).Value = '  +10.87%  '

# Row 30
Set-TextValue 30 4 '44.01'
$ws.Cells.Item(30, 5).Value = '  +6.94%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.63%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -0.35%  '

# Row 33
Set-TextValue 33 4 '11.84'
$ws.Cells.Item(33, 5).Value = '  +5.08%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -0.13%  '

# Row 35
Set-TextValue 35 4 '0.0493'
$ws.Cells.Item(35, 5).Value = '  +2.84%  '

# Row 36
Set-TextValue 36 4 '52.21'
$ws.Cells.Item(36, 5).Value = '  +0.70%  '

# Row 37
Set-TextValue 37 4 '0.998'
$ws.Cells.Item(37, 5).Value = '  +0.24%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.98%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +0.56%  '

# Row 40
Set-TextValue 40 4 '0.132'
$ws.Cells.Item(40, 5).Value = '  +6.85%  '

# Row 41
Set-TextValue 41 4 '144.33'
$ws.Cells.Item(41, 5).Value = '  +5.38%  '

# Row 42
Set-TextValue 42 4 '0.313'
$ws.Cells.Item(42, 5).Value = '  +9.42%  '

# Row 43
Set-TextValue 43 4 '1.97'
$ws.Cells.Item(43, 5).Value = '  -0.82%  '

# Row 44
Set-TextValue 44 4 '4.00'
$ws.Cells.Item(44, 5).Value = '  -0.63%  '

# Row 45
Set-TextValue 45 4 '16.72'
$ws.Cells.Item(45, 5).Value = '  +0.53%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +4.54%  '

# Row 47
Set-TextValue 47 4 '21.76'
$ws.Cells.Item(47, 5).Value = '  +2.00%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -0.42%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -0.36%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +2.84%  '

# Row 51
Set-TextValue 51 4 '0.0363'
$ws.Cells.Item(51, 5).Value = '  +6.07%  '
